# 16.4.2.1 — add a "(units)" sub-header row and extend the table with 2022/2023 columns.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new row under the title row for the "(units)" / "(единиц)" /
#    "(бирдик)" annotation, shifting the old rows 2-5 down to 3-6.
# ---------------------------------------------------------------------------
$ws.Rows(2).EntireRow.Insert()

$ws.Range("A2").Value = "(бирдик)"
$ws.Range("B2").Value = "(единиц)"
$ws.Range("C2").Value = "(units)"

$ws.Range("A2:C2").Font.Name = "Times New Roman"
$ws.Range("A2:C2").Font.Size = 9
$ws.Range("A2:C2").Font.Italic = $true
$ws.Range("A2:C2").HorizontalAlignment = -4108
$ws.Range("A2:C2").VerticalAlignment = -4108

# ---------------------------------------------------------------------------
# 2. Add the 2022 / 2023 columns (J, K) to the header and data rows.
# ---------------------------------------------------------------------------
$ws.Range("J4").Value = 2022
$ws.Range("K4").Value = 2023
$ws.Range("J4:K4").Style = $ws.Range("I4").Style

$ws.Range("J5").Value = "-"
$ws.Range("K5").Value = 219
$ws.Range("J5").Style = $ws.Range("I5").Style
$ws.Range("K5").Style = $ws.Range("I5").Style

$ws.Range("J6").Value = "-"
$ws.Range("K6").Value = 171
$ws.Range("J6").Style = $ws.Range("I6").Style
$ws.Range("K6").Style = $ws.Range("I6").Style

$ws.Range("J5:J6").HorizontalAlignment = -4152
$ws.Range("J5:J6").Font.Name = "Times New Roman"
$ws.Range("J5:J6").Font.Size = 9

# ---------------------------------------------------------------------------
# 3. Cosmetic tweaks to match the refreshed layout.
# ---------------------------------------------------------------------------
$ws.Rows(1).RowHeight = 41.25
$ws.Rows(6).RowHeight = 18

$ws.Columns("A:C").ColumnWidth = 34.99
